$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe the old 5-player table (rows 2-6, cols A-J) entirely.
$ws.Cells.Clear()

# New header row: Equipo / Jugador replace the old Nombre column, rest shifts right one column.
$ws.Range("A1").Value = "Equipo"
$ws.Range("B1").Value = "Jugador"
$ws.Range("C1").Value = "TCA"
$ws.Range("D1").Value = "T3A"
$ws.Range("E1").Value = "TCI"
$ws.Range("F1").Value = "TLI"
$ws.Range("G1").Value = "Puntos"
$ws.Range("H1").Value = "%FG"
$ws.Range("I1").Value = "%eFG"
$ws.Range("J1").Value = "%TS"
$ws.Range("K1").Value = "Valoration"

# Row 2 - Angeles Lakers / Anthony Davis
$ws.Range("A2").Value = "Angeles Lakers"
$ws.Range("B2").Value = "Anthony Davis"
$ws.Range("C2").Value = 4.0
$ws.Range("D2").Value = 2.0
$ws.Range("E2").Value = 6.0
$ws.Range("F2").Value = 2.0
$ws.Range("G2").Value = 12.0
$ws.Range("H2").Value = 66.67
$ws.Range("I2").Value = 83.33
$ws.Range("J2").Value = 87.21
$ws.Range("K2").Value = 24.0

# Row 3 - Golden State Warriors / Klay Thompson
$ws.Range("A3").Value = "Golden State Warriors"
$ws.Range("B3").Value = "Klay Thompson"
$ws.Range("C3").Value = 4.0
$ws.Range("D3").Value = 2.0
$ws.Range("E3").Value = 6.0
$ws.Range("F3").Value = 234.0
$ws.Range("G3").Value = 244.0
$ws.Range("H3").Value = 66.67
$ws.Range("I3").Value = 83.33
$ws.Range("J3").Value = 111.97
$ws.Range("K3").Value = 277.0

# Row 4 - promedios (averages of the two data rows above)
$ws.Range("A4").Value = "promedios"
$ws.Range("B4").Formula = "=AVERAGE(B2:B3)"
$ws.Range("C4").Formula = "=AVERAGE(C2:C3)"
$ws.Range("D4").Formula = "=AVERAGE(D2:D3)"
$ws.Range("E4").Formula = "=AVERAGE(E2:E3)"
$ws.Range("F4").Formula = "=AVERAGE(F2:F3)"
$ws.Range("G4").Formula = "=AVERAGE(G2:G3)"
$ws.Range("H4").Formula = "=AVERAGE(H2:H3)"
$ws.Range("I4").Formula = "=AVERAGE(I2:I3)"
$ws.Range("J4").Formula = "=AVERAGE(J2:J3)"
$ws.Range("K4").Formula = "=AVERAGE(K2:K3)"

# Match the saved selection state (D3:D4, active cell D4).
$ws.Range("D3:D4").Select()
